# "fixes in reconnection" - append 3 new reconnection-event rows (for the
# three new timestamps) to each of the three message sheets.

$wb = $excel.ActiveWorkbook

$ts1 = "2022-07-02 13:57:07"
$ts2 = "2022-07-03 05:03:00"
$ts3 = "2022-07-03 08:21:19"

# ---------------------------------------------------------------------
# Sheet "Msg8705": rows 187-189, columns A:H
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Msg8705")

$ws1.Cells.Item(187, 1).Value = $ts1
$ws1.Cells.Item(187, 2).Value = ""
$ws1.Cells.Item(187, 3).Value = "10F872226797"
$ws1.Cells.Item(187, 4).Value = 1656763028
$ws1.Cells.Item(187, 5).Value = 3
$ws1.Cells.Item(187, 6).Value = 12
$ws1.Cells.Item(187, 7).Value = 1
$ws1.Cells.Item(187, 8).Value = 255

$ws1.Cells.Item(188, 1).Value = $ts2
$ws1.Cells.Item(188, 2).Value = ""
$ws1.Cells.Item(188, 3).Value = "10F872226797"
$ws1.Cells.Item(188, 4).Value = 1656817385
$ws1.Cells.Item(188, 5).Value = 3
$ws1.Cells.Item(188, 6).Value = 12
$ws1.Cells.Item(188, 7).Value = 1
$ws1.Cells.Item(188, 8).Value = 255

$ws1.Cells.Item(189, 1).Value = $ts3
$ws1.Cells.Item(189, 2).Value = ""
$ws1.Cells.Item(189, 3).Value = "10F872226797"
$ws1.Cells.Item(189, 4).Value = 1656829284
$ws1.Cells.Item(189, 5).Value = 3
$ws1.Cells.Item(189, 6).Value = 12
$ws1.Cells.Item(189, 7).Value = 1
$ws1.Cells.Item(189, 8).Value = 255

# ---------------------------------------------------------------------
# Sheet "Msg8705_8": rows 187-189, columns A:L
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Msg8705_8")

$ws2.Cells.Item(187, 1).Value = $ts1
$ws2.Cells.Item(187, 2).Value = ""
$ws2.Cells.Item(187, 3).Value = 522
$ws2.Cells.Item(187, 4).Value = 37122
$ws2.Cells.Item(187, 5).Value = 655618
$ws2.Cells.Item(187, 6).Value = 10485764
$ws2.Cells.Item(187, 7).Value = 1
$ws2.Cells.Item(187, 8).Value = ""
$ws2.Cells.Item(187, 9).Value = "1.0.1"
$ws2.Cells.Item(187, 10).Value = ""
$ws2.Cells.Item(187, 11).Value = 21
$ws2.Cells.Item(187, 12).Value = ""

$ws2.Cells.Item(188, 1).Value = $ts2
$ws2.Cells.Item(188, 2).Value = ""
$ws2.Cells.Item(188, 3).Value = 522
$ws2.Cells.Item(188, 4).Value = 37122
$ws2.Cells.Item(188, 5).Value = 655618
$ws2.Cells.Item(188, 6).Value = 10485764
$ws2.Cells.Item(188, 7).Value = 1
$ws2.Cells.Item(188, 8).Value = ""
$ws2.Cells.Item(188, 9).Value = "1.0.1"
$ws2.Cells.Item(188, 10).Value = ""
$ws2.Cells.Item(188, 11).Value = 21
$ws2.Cells.Item(188, 12).Value = ""

$ws2.Cells.Item(189, 1).Value = $ts3
$ws2.Cells.Item(189, 2).Value = ""
$ws2.Cells.Item(189, 3).Value = 522
$ws2.Cells.Item(189, 4).Value = 37122
$ws2.Cells.Item(189, 5).Value = 0
$ws2.Cells.Item(189, 6).Value = 0
$ws2.Cells.Item(189, 7).Value = 1
$ws2.Cells.Item(189, 8).Value = ""
$ws2.Cells.Item(189, 9).Value = "1.0.1"
$ws2.Cells.Item(189, 10).Value = ""
$ws2.Cells.Item(189, 11).Value = 21
$ws2.Cells.Item(189, 12).Value = ""

# ---------------------------------------------------------------------
# Sheet "Msg8705_11": rows 557-565, columns A:H
# Each of the three new timestamps gets three rows (one per C-value).
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Msg8705_11")

$cVals = @(1914729089, 1914726672, 1914728742)
$timestamps = @($ts1, $ts2, $ts3)

$r = 557
foreach ($ts in $timestamps) {
    foreach ($cVal in $cVals) {
        $ws3.Cells.Item($r, 1).Value = $ts
        $ws3.Cells.Item($r, 2).Value = ""
        $ws3.Cells.Item($r, 3).Value = $cVal
        $ws3.Cells.Item($r, 4).Value = 10012
        $ws3.Cells.Item($r, 5).Value = 269627400
        $ws3.Cells.Item($r, 6).Value = 256
        $ws3.Cells.Item($r, 7).Value = 2560
        $ws3.Cells.Item($r, 8).Value = 8193
        $r = $r + 1
    }
}
